$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Dataset Registry"
# Replace the old single DS002 (sep6p6) row with five new rows
# (DS001-DS005) of sep5p1 data, seeds 4,0,1,3,2 respectively.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Dataset Registry")

$ds1Rows = @(
    @("DS001", "n1000000_f_init5_cont0_disc5_sep5p1_seed4_config.yml", "n1000000_f_init5_cont0_disc5_sep5p1_seed4_dataset.csv", "2025-07-27", "1,000,000 samples, 5 features, Avg. Sep: 5.10"),
    @("DS002", "n1000000_f_init5_cont0_disc5_sep5p1_seed0_config.yml", "n1000000_f_init5_cont0_disc5_sep5p1_seed0_dataset.csv", "2025-07-27", "1,000,000 samples, 5 features, Avg. Sep: 5.10"),
    @("DS003", "n1000000_f_init5_cont0_disc5_sep5p1_seed1_config.yml", "n1000000_f_init5_cont0_disc5_sep5p1_seed1_dataset.csv", "2025-07-27", "1,000,000 samples, 5 features, Avg. Sep: 5.10"),
    @("DS004", "n1000000_f_init5_cont0_disc5_sep5p1_seed3_config.yml", "n1000000_f_init5_cont0_disc5_sep5p1_seed3_dataset.csv", "2025-07-27", "1,000,000 samples, 5 features, Avg. Sep: 5.10"),
    @("DS005", "n1000000_f_init5_cont0_disc5_sep5p1_seed2_config.yml", "n1000000_f_init5_cont0_disc5_sep5p1_seed2_dataset.csv", "2025-07-27", "1,000,000 samples, 5 features, Avg. Sep: 5.10")
)

for ($i = 0; $i -lt $ds1Rows.Count; $i++) {
    $r = $i + 2
    $row = $ds1Rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]

    # The "Creation Date" column holds a plain text value that looks like a
    # date (e.g. "2025-07-27"). Force it to stay text (not get silently
    # converted into a date serial number) by writing it with a text
    # number format, then reset the cell style back to Normal so no
    # residual formatting is left attached to the cell.
    $dcell = $ws1.Cells.Item($r, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $row[3]
    $dcell.Style = "Normal"

    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------
# Sheet 2: "Configuration Details"
# Remove the old DS002 data row, leaving only the header row.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Configuration Details")
$ws2.Range("A2:H2").Clear()

# ---------------------------------------------------------------
# Sheet 3: "Feature Separation Details"
# Remove the old DS002 feature rows (2-6), leaving only the header.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Feature Separation Details")
$ws3.Range("A2:H6").Clear()

# ---------------------------------------------------------------
# Sheet 4: "File Metadata"
# Remove the old DS002 data row, leaving only the header row.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("File Metadata")
$ws4.Range("A2:E2").Clear()
